# "cambios de may de mayo" -- update the UPP Pachuca quejas/sugerencias report
# from the 4th-quarter-2021 reporting period to the 1st-quarter-2022 period,
# refresh the validation/update dates, and replace one incident narrative with
# a new one (plus small formatting touch-ups picked up along the way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header block (row 3): wrap the "DESCRIPCIÓN" merged header and its blanks
# ---------------------------------------------------------------------------
$ws.Rows(3).RowHeight = 29.25
$ws.Range("G3").WrapText = $true
$ws.Range("H3:I3").WrapText = $true

# Only a left border remains on the merged DESCRIPCIÓN header cell
$ws.Range("G3").Borders(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("G3").Borders(10).LineStyle = -4142  # xlEdgeRight  -> none
$ws.Range("G3").Borders(8).LineStyle = -4142   # xlEdgeTop    -> none

# ---------------------------------------------------------------------------
# Column width tweaks
# ---------------------------------------------------------------------------
$ws.Columns(7).ColumnWidth = 81.67   # column G ("Resultados")
$ws.Columns(11).ColumnWidth = 21.83  # column K ("Nota")

# ---------------------------------------------------------------------------
# Data rows 8-11: move the reporting period from Q4 2021 to Q1 2022 and push
# the validation / update dates from 2022-02-10 to 2022-04-08
# ---------------------------------------------------------------------------
foreach ($r in 8..11) {
    $ws.Cells.Item($r, 1).Value = 2022     # A: Ejercicio
    $ws.Cells.Item($r, 2).Value = 44562    # B: Fecha de inicio (2022-01-01)
    $ws.Cells.Item($r, 3).Value = 44651    # C: Fecha de término (2022-03-31)
    $ws.Cells.Item($r, 9).Value = 44659    # I: Fecha de validación (2022-04-08)
    $ws.Cells.Item($r, 10).Value = 44659   # J: Fecha de Actualización (2022-04-08)
}

# Row 8 ("Facebook")
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Rows(8).RowHeight = 111

# Row 9 ("Correo electrónico") -- new complaint narrative, participant count
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "Un alumno que esta realizando sus trámites de Estadia no tenia aun firmadas unas calificaciones y esto no le permitia continuar con su trámite. El área académica responde que de acuerdo con el reglamento de licenciatura de esta casa de estudios, en el artículo 90 se expresa el tiempo de publicación, revisión y firma de los resultados en el Sistema Integral de la Universidad Politécnica de Pachuca aun se estaba en tiempo y se firmarón las calificaciones."
$ws.Range("K9").Value = "Los campos que se observan vacío es por que no se presentaron quejas por este medio"
$ws.Rows(9).RowHeight = 78

# Row 10 ("Buzón de quejas y sugerencias")
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Rows(10).RowHeight = 107.25

# Row 11 ("Llamada telefónica")
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Rows(11).RowHeight = 94.5

# ---------------------------------------------------------------------------
# View state: scrolled down to row 11, cursor left on C19
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("C19").Select()
